$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7848485112190247
$ws.Range("B1").Value = 1.147804141044617
$ws.Range("C1").Value = 3.570649147033691
$ws.Range("D1").Value = 3.846026182174683
$ws.Range("E1").Value = 1.867638230323792
